# Doing Updates for Financials
#
# Refresh the yearly financial figures on the "HMC" sheet (columns D..J
# hold the 7 most-recent fiscal years, newest first) with updated source
# data. Only the numeric figures change; row/column layout, labels, and
# formatting are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("HMC")

# Map of row number -> new values for columns D..J (columns 4..10).
# $null entries mean "leave the existing cell value/content alone"
# (e.g. row 49 keeps its "NA" text in columns I/J).
$rowData = @{
    8   = @(138864800, 126552800, 131994400, 120486000, 113055100, 89296600, 71850800)
    9   = @(108485300, 98227300, 102444900, 93390300, 86698600, 66400300, 53513500)
    10  = @(30379500, 28325500, 29549500, 27095700, 26356400, 22896400, 18337300)
    12  = @(6796800, 6250500, 5934800, 5479700, 5409300, 5064800, 4699200)
    17  = @(131329400, 118952700, 127443900, 114423800, 105607300, 84371600, 69759200)
    18  = @(7535400, 7600000, 4550500, 6062300, 7447700, 4925100, 2091500)
    20  = @(2661200, 1615900, 1358000, 1390600, 1110500, -395600, 329200)
    21  = @(16650100, 15318600, 11888000, 13111200, 13880800, 9873200, 7442300)
    22  = @(117200, 112700, 164000, 164500, 115700, 109900, 93800)
    23  = @(10079400, 9103200, 5744500, 7288400, 8442500, 4419600, 2326900)
    24  = @(-123500, 2961400, 2071000, 2216100, 2422600, 1617900, 1227000)
    26  = @(10202900, 6141700, 3673500, 5072300, 6019800, 2801600, 1099900)
    27  = @(9576400, 5573800, 3114600, 4605300, 6844900, 3319000, 1911800)
    32  = @(-2661200, -1615900, -1358000, -1390600, -1110500, 395600, -329200)
    33  = @(9576400, 5573800, 3114600, 4605300, 6844900, 3319000, 1911800)
    35  = @(9576400, 5573800, 3114600, 4605300, 6844900, 3319000, 1911800)
    41  = @(20398700, 19038000, 15887400, 13304400, 10790000, 10903400, 22547800)
    42  = @(1927100, 1350800, 931400, 838100, 922300, 3782800, 3377000)
    43  = @(23876100, 23892400, 24884700, 26393500, 34629400, 20330800, 34241300)
    44  = @(13772000, 12331700, 11872200, 13544700, 16351900, 10987400, 18726900)
    45  = @(2630700, 2648400, 2848600, 2836400, 4058000, 2116000, 6789700)
    46  = @(62604600, 59261400, 56424300, 56917100, 50164400, 48120400, 42841300)
    47  = @(38270300, 36453700, 36252700, 41133900, 43665200, 29355100, 50608400)
    48  = @(64641100, 66037600, 61631800, 58984900, 93149100, 38353700, 65025200)
    49  = @(6703300, 7034900, 7457400, 6866200, 6054800, $null, $null)
    52  = @(2697200, 2593900, 3026500, 2667500, 11104600, 7434500, 8759700)
    54  = @(174916000, 171381000, 164793000, 166570000, 145078000, 123264000, 106560000)
    57  = @(9533500, 9160300, 8692900, 9036300, 18126400, 8931600, 17518500)
    58  = @(26372000, 25193800, 25218200, 25615400, 47417000, 19737400, 33922500)
    59  = @(14936300, 14722600, 15540900, 13269900, 11453600, 8295500, 13281100)
    60  = @(50841900, 49076700, 49452000, 47921500, 42956300, 36964600, 32361000)
    61  = @(35633500, 36787700, 34210800, 36046300, 29809400, 24506000, 20204400)
    62  = @(14004900, 17087600, 17562700, 15861000, 27151100, 14736000, 26149500)
    66  = @(103197000, 105432000, 103670000, 102308000, 87804700, 77670400, 66854200)
    72  = @(68806400, 60684600, 55996600, 54995500, 111301000, 54685100, 104637000)
    76  = @(71719200, 65949500, 61123400, 64262000, 57273200, 45593200, 39705700)
    81  = @(9576400, 5573800, 3114600, 4605300, 6844900, 3319000, 1911800)
    83  = @(6446400, 6095900, 5972900, 5652100, 5316700, 5337800, 5016000)
    89  = @(8928500, 8001100, 12574600, 9224500, 4059900, 7238700, 6884300)
    91  = @(-3756700, -4467000, -5742000, -5859800, -23687300, -12836800, -9772100)
    94  = @(-5560600, -5881600, -7910700, -7598100, -8326000, -9670600, -6084500)
    96  = @(-1575000, -1466300, -1433800, -1433800, -1287100, -1173100, -977600)
    100 = @(-1576000, 1043400, -861500, 112900, 3213300, 1080900, -616800)
    101 = @(-431300, -12300, -1219400, 775200, 356400, 980500, -471400)
    102 = @(1360600, 3150600, 2583000, 2514400, -696300, -370500, -288500)
}

foreach ($row in $rowData.Keys) {
    $vals = $rowData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        if ($null -ne $vals[$i]) {
            $ws.Cells.Item([int]$row, 4 + $i).Value = $vals[$i]
        }
    }
}
